$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (Mardi, 2nd table): add Fin time, set Temps total, update Travail effectué text
$ws.Range("D15").Value = 0.64583333333333337
$ws.Range("D15").NumberFormat = "h:mm"
$ws.Range("E15").Value = "7h00"
$ws.Range("F15").Value = "Agenda, modification d'éléments mineurs détectés, devis double clic"

# Row 16 (Mercredi, 2nd table): add Début time, set Travail effectué text
$ws.Range("B16").Value = 0.3527777777777778
$ws.Range("B16").NumberFormat = "h:mm"
$ws.Range("F16").Value = "devis double clic"

# Update the active selection to match the saved view state
$ws.Range("F17").Select()
